$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Copy()
$ws.Range("B338:M338").PasteSpecial(-4122)
